$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.617
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.598000000000001
$ws.Range("B25").Value = 6.425999999999999
$ws.Range("D27").Value = -8.625999999999999
$ws.Range("D31").Value = -8.186
$ws.Range("D39").Value = -8.063000000000001
$ws.Range("D48").Value = -7.31
$ws.Range("D51").Value = -8.34
$ws.Range("D52").Value = -8.100000000000001
$ws.Range("B53").Value = 6.076
$ws.Range("D55").Value = -7.904000000000001
$ws.Range("D56").Value = -8.434999999999999
$ws.Range("B57").Value = 4.981999999999999
$ws.Range("D57").Value = -8.059999999999999
$ws.Range("B59").Value = 5.145
$ws.Range("B69").Value = 5.339
$ws.Range("D73").Value = -8.278000000000002
$ws.Range("B79").Value = 5.411
$ws.Range("B83").Value = 5.915
$ws.Range("D89").Value = -6.173999999999999
$ws.Range("D90").Value = -7.49
$ws.Range("B93").Value = 5.608
